$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 4, shifting existing rows 4-25 down to 5-26.
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the latest weekly price record.
$ws.Cells.Item(4, 1).Value = 1
$ws.Cells.Item(4, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(4, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(4, 4).Value = 45251
$ws.Cells.Item(4, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(4, 5).Value = 15
$ws.Cells.Item(4, 6).Value = "Fruta"
$ws.Cells.Item(4, 7).Value = 100103
$ws.Cells.Item(4, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(4, 9).Value = 100103001
$ws.Cells.Item(4, 10).Value = "Cereza"
$ws.Cells.Item(4, 11).Value = "Early Burlat"
$ws.Cells.Item(4, 12).Value = "Primera"
$ws.Cells.Item(4, 13).Value = 300
$ws.Cells.Item(4, 14).Value = 14000
$ws.Cells.Item(4, 15).Value = 15000
$ws.Cells.Item(4, 16).Value = 14500
$ws.Cells.Item(4, 17).Value = "`$/bandeja 5 kilos"
$ws.Cells.Item(4, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(4, 19).Value = 2900
$ws.Cells.Item(4, 20).Value = 5
